$d = $word.ActiveDocument

# 1. Remove the "_GoBack" bookmark from its old location (after "[date]")
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

# 2. Remove "Dra. " from the signature line ("Dra. Laura Ación" -> "Laura Ación")
$d.Content.Find.Execute("Dra. Laura Ación", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Laura Ación", 2)

# 3. Re-create the "_GoBack" bookmark at the start of the edited run, marking
#    the location of the last edit (mirrors Word's own behaviour).
$rng = $d.Content
$rng.Find.Execute("Laura Ación", $true, $false, $false, $false, $false,
                   $true, 1, $false, "", 0)
$newBookmarkRange = $d.Range($rng.Start, $rng.Start)
$d.Bookmarks.Add("_GoBack", $newBookmarkRange)
